# RO dist. small cities (2018)
# Adds 4 new small Rostov-region cities (Гуково, Донецк, Каменск-Шахтинский,
# Зверево) for year 2018, fills in missing "???" placeholders for Азов 2022
# (row 35, columns G:K), and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# Row 35 (Азов, 2022): fill previously-empty G:K with the "???" marker,
# matching the formatting already used on the rest of the row (e.g. P35).
# ---------------------------------------------------------------------
$rng = $ws.Range("G35:K35")
$rng.Value = "???"
$rng.HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Helper data for the four new rows (all year 2018).
# ---------------------------------------------------------------------
$rowsData = @(
    @{
        Row = 38
        Name = "Гуково"
        C = 64.869
        D = 8.482
        E = 958
        F = 22958.3
        OFormula = "=774396/1000"
        Q = $null
        QFormula = "=8137/1000"
        SFormula = "=6362296.4/1000"
        TFormula = "=221935.5/1000"
        U = -403
    },
    @{
        Row = 39
        Name = "Донецк"
        C = 47770
        D = 6.41
        E = 523
        F = 21351
        OFormula = "=2584304/1000"
        Q = 7.315
        QFormula = $null
        SFormula = "=4911240.3/1000"
        TFormula = "=138739.6/1000"
        U = -197
    },
    @{
        Row = 40
        Name = "Каменск-Шахтинский"
        C = 88.997
        D = 27.875
        E = 753
        F = 28590.4
        OFormula = "=20336870/1000"
        Q = 16.305
        QFormula = $null
        SFormula = "=17686635.8/1000"
        TFormula = "=742584.8/1000"
        U = -71
    },
    @{
        Row = 41
        Name = "Зверево"
        C = 19.045
        D = 5.342
        E = 382
        F = 27277.2
        O = 209.603
        OFormula = $null
        Q = 2.024
        QFormula = $null
        SFormula = "=1957826.5/1000"
        TFormula = "=53897.4/1000"
        U = -83
    }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # A:B -- city name / year, bold + centered (matches the rest of the table)
    $ws.Range("A$r").Value = $rd.Name
    $ws.Range("B$r").Value = 2018
    $ab = $ws.Range("A$r`:B$r")
    $ab.HorizontalAlignment = $xlCenter
    $ab.Font.Bold = $true

    # C:F -- plain numeric stats, centered
    $ws.Range("C$r").Value = $rd.C
    $ws.Range("D$r").Value = $rd.D
    $ws.Range("E$r").Value = $rd.E
    $ws.Range("F$r").Value = $rd.F
    $ws.Range("C$r`:F$r").HorizontalAlignment = $xlCenter

    # G:N -- unavailable data, "???" marker, centered
    $gn = $ws.Range("G$r`:N$r")
    $gn.Value = "???"
    $gn.HorizontalAlignment = $xlCenter

    # O -- invests (thousands), formula or literal depending on row
    if ($rd.OFormula) {
        $ws.Range("O$r").Formula = $rd.OFormula
    } else {
        $ws.Range("O$r").Value = $rd.O
    }
    $ws.Range("O$r").HorizontalAlignment = $xlCenter

    # P -- unavailable data, "???" marker, centered
    $p = $ws.Range("P$r")
    $p.Value = "???"
    $p.HorizontalAlignment = $xlCenter

    # Q -- companies, formula or literal depending on row
    if ($rd.QFormula) {
        $ws.Range("Q$r").Formula = $rd.QFormula
    } else {
        $ws.Range("Q$r").Value = $rd.Q
    }
    $ws.Range("Q$r").HorizontalAlignment = $xlCenter

    # R -- unavailable data, "???" marker, centered
    $rr = $ws.Range("R$r")
    $rr.Value = "???"
    $rr.HorizontalAlignment = $xlCenter

    # S:T -- retail / foodserv turnover (thousands), formulas
    $ws.Range("S$r").Formula = $rd.SFormula
    $ws.Range("T$r").Formula = $rd.TFormula
    $ws.Range("S$r`:T$r").HorizontalAlignment = $xlCenter

    # U -- saldo, centered
    $ws.Range("U$r").Value = $rd.U
    $ws.Range("U$r").HorizontalAlignment = $xlCenter
}

# ---------------------------------------------------------------------
# View: drop the frozen top-left scroll position and move the selection.
# ---------------------------------------------------------------------
$ws.Range("W25").Select()
